$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: remove old row 1 content (C1, E1) and clear all data rows except row 2 (header, untouched to avoid row-height autofit)
$ws.Range("A1:F1").ClearContents()
$ws.Range("A3:F29").ClearContents()

# Step 2: write the new cell contents
$ws.Range("A3").Value = "Fall 2022"
$ws.Range("B3").Value = "Credits"
$ws.Range("C3").Value = "Spring 2022"
$ws.Range("D3").Value = "Credits"
$ws.Range("E3").Value = "Summer 2022"
$ws.Range("F3").Value = "Credits"
$ws.Range("A4").Value = "POLS 1101"
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = "CPSC 3165"
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = "CPSC 4698"
$ws.Range("F4").Value = 3
$ws.Range("A5").Value = "DSCI 3111"
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = "CPSC 3415"
$ws.Range("D5").Value = 1
$ws.Range("A6").Value = "ARTH 3115"
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = "CPSC 4135"
$ws.Range("D6").Value = 3
$ws.Range("A7").Value = "ARTH 3119"
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = "CPSC 4148"
$ws.Range("D7").Value = 3
$ws.Range("A8").Value = "CPSC 3121"
$ws.Range("B8").Value = 3
$ws.Range("C8").Value = "CPSC 4155"
$ws.Range("D8").Value = 3
$ws.Range("A9").Value = "CPSC 4000"
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = "CYBR 4416"
$ws.Range("D9").Value = 1
$ws.Range("A11").Value = "Total"
$ws.Range("B11").Formula = "=SUM(B4:B10)"
$ws.Range("C11").Value = "Total"
$ws.Range("D11").Formula = "=SUM(D4:D10)"
$ws.Range("E11").Value = "Total"
$ws.Range("F11").Formula = "=SUM(F4:F10)"
$ws.Range("A12").Value = "Fall 2023"
$ws.Range("B12").Value = "Credits"
$ws.Range("C12").Value = "Spring 2023"
$ws.Range("D12").Value = "Credits"
$ws.Range("E12").Value = "Summer 2023"
$ws.Range("F12").Value = "Credits"
$ws.Range("A13").Value = "CPSC 4157"
$ws.Range("B13").Value = 3
$ws.Range("C13").Value = "CPSC 4176"
$ws.Range("D13").Value = 3
$ws.Range("A14").Value = "CPSC 4175"
$ws.Range("B14").Value = 3
$ws.Range("A15").Value = "CPSC 4205"
$ws.Range("B15").Value = 3
$ws.Range("A20").Value = "Total"
$ws.Range("B20").Formula = "=SUM(B13:B19)"
$ws.Range("C20").Value = "Total"
$ws.Range("D20").Formula = "=SUM(D13:D19)"
$ws.Range("E20").Value = "Total"
$ws.Range("F20").Formula = "=SUM(F13:F19)"
$ws.Range("A21").Value = "Fall 2024"
$ws.Range("B21").Value = "Credits"
$ws.Range("C21").Value = "Spring 2024"
$ws.Range("D21").Value = "Credits"
$ws.Range("E21").Value = "Summer 2024"
$ws.Range("F21").Value = "Credits"
$ws.Range("A29").Value = "Total"
$ws.Range("B29").Formula = "=SUM(B22:B28)"
$ws.Range("C29").Value = "Total"
$ws.Range("D29").Formula = "=SUM(D22:D28)"
$ws.Range("E29").Value = "Total"
$ws.Range("F29").Formula = "=SUM(F22:F28)"
